# InstallTracker.xlsx update
# - Adds 4 more days of "Actual" install counts (and corrects the most
#   recent day's count), which ripples through the weekly growth-rate
#   column and re-fits the Linear / Poly-2 / Poly-3 trend coefficients.
# - Moves the remembered cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Actual install counts (column C) ------------------------------
# Row 60 (day 58) had a later correction; rows 61-64 (days 59-62) are
# brand-new data points that previously had no "Actual" value at all.
$ws.Range("C60").Value = 593
$ws.Range("C61").Value = 604
$ws.Range("C62").Value = 647
$ws.Range("C63").Value = 687
$ws.Range("C64").Value = 729

# --- 2. Growth/Day from Last Week (column D) ---------------------------
# D2 is the (literal) first-week value; with the corrected history it
# becomes 0.  D61:D64 are newly populated with the same "(Cn-C(n-7))/7"
# pattern used by the rest of the column.
$ws.Range("D2").Value = 0

$ws.Range("D61").Formula = "=(C61-C54)/7"
$ws.Range("D62").Formula = "=(C62-C55)/7"
$ws.Range("D63").Formula = "=(C63-C56)/7"
$ws.Range("D64").Formula = "=(C64-C57)/7"

# --- 3. Re-fit trend formulas (columns E/F/G) ---------------------------
# Row 2 holds the "anchor" (non-shared) copy of each formula; the bulk of
# each column is one shared formula for rows 3-66 and another for rows
# 67-115 (row 66/67 is where Excel originally split the fill).
$ws.Range("E2").Formula = "=7.8126*B2+3"
$ws.Range("F2").Formula = "=0.2073*B2^2-1.9062*B2+3"
$ws.Range("G2").Formula = "=0.0003*B2^3+0.1854*B2^2-1.4945*B2+3"

$ws.Range("E3:E66").Formula = "=7.8126*B3+3"
$ws.Range("F3:F66").Formula = "=0.2073*B3^2-1.9062*B3+3"
$ws.Range("G3:G66").Formula = "=0.0003*B3^3+0.1854*B3^2-1.4945*B3+3"

$ws.Range("E67:E115").Formula = "=7.8126*B67+3"
$ws.Range("F67:F115").Formula = "=0.2073*B67^2-1.9062*B67+3"
$ws.Range("G67:G115").Formula = "=0.0003*B67^3+0.1854*B67^2-1.4945*B67+3"

# --- 4. Remembered selection --------------------------------------------
$null = $ws.Range("I20").Select()
